$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, station, fee income, total income, charge amount, order count)
$rows = @(
    @(46011, "四方坪站", 7786.07, 6608.19, 2513.11, 352),
    @(46011, "高岭站",   4371.05, 3812.57, 1163.87, 156),
    @(46012, "四方坪站", 9150.84, 7944.04, 2926.25, 385),
    @(46012, "高岭站",   4651.65, 3844.94, 1228.97, 157)
)

$startRow = 40
$lastExistingRow = 39

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Clone formatting (and base values) from the last existing data row so the
    # new rows reuse the existing cell styles instead of creating new ones.
    $ws.Range("A" + $lastExistingRow + ":F" + $lastExistingRow).Copy()
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4104)

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

$excel.CutCopyMode = $false

$ws.Application.GoTo($ws.Range("A34"))
$ws.Range("I43").Select()
